$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in grades (5) for row 24 (Слепцов Илья) columns C:E
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 5

# Fill in grades (5) for row 25 (Слепцов Лев) columns D:E (C25 already had 5)
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 5

# Move the active selection to E26, matching the author's final cursor position
$ws.Range("E26").Select()
